$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Semana 36 de 2025: remove event "610" (row without esperado/observado/valor-p data)
# and event "760" Tetanos accidental (discontinued), then update weekly counts.
$ws.Rows("33").Delete()
$ws.Rows("29").Delete()

# Update Esperado (C), Observado (D) and valor p (E) for remaining events
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 0.1
$ws.Range("D4").Value = 0
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 0.13
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 7
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 0.27
$ws.Range("C9").Value = 43
$ws.Range("D9").Value = 38
$ws.Range("E9").Value = 0.05
$ws.Range("C10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 0.27
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 0.18
$ws.Range("C13").Value = 3
$ws.Range("E13").Value = 0.05
$ws.Range("C15").Value = 2
$ws.Range("E15").Value = 0.14
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 1
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = 0.04
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 0.14
$ws.Range("C19").Value = 11
$ws.Range("E19").Value = 0.04
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 1
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 0.15
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = 0.37
$ws.Range("C29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("C31").Value = 3
$ws.Range("E31").Value = 0.22
$ws.Range("C33").Value = 9
$ws.Range("D33").Value = 6
$ws.Range("E33").Value = 0.09
$ws.Range("C34").Value = 9
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 0
$ws.Range("C35").Value = 7
$ws.Range("D35").Value = 5
$ws.Range("E35").Value = 0.13
